$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.868.86'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.809.96'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.55'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4623'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3698'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07351'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8749'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.45'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.813.77'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.358'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.506'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.73'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07043'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008694'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.857.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.313'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.957.89'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.897'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.59'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.37'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.149'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.315'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.95'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08897'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7531'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.154'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.30%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.920'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.67%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.458'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.100'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01963'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.443'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05250'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.921'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5311'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.148'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1659'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.452'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4953'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.36'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.671'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.35'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06287'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.53%  '

Write-Host "Applied changes"